$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-01 Thursday" "2026-01-02 Friday"

Replace-Text "39×97=" "94×31="
Replace-Text "31×49=" "39×24="
Replace-Text "37×82=" "57×41="
Replace-Text "58×95=" "27×82="
Replace-Text "97×76=" "80×20="

Replace-Text "60×70=" "23×97="
Replace-Text "84×35=" "95×89="
Replace-Text "79×15=" "38×77="
Replace-Text "90×58=" "18×84="
Replace-Text "72×71=" "79×91="

Replace-Text "87×13=" "68×28="
Replace-Text "55×98=" "65×38="
Replace-Text "92×74=" "34×95="
Replace-Text "57×86=" "68×93="
Replace-Text "49×55=" "63×60="

Replace-Text "12×78=" "63×17="
Replace-Text "55×73=" "44×41="
Replace-Text "86×59=" "95×52="
Replace-Text "30×63=" "53×69="
Replace-Text "50×41=" "34×86="

Replace-Text "53×76=" "53×73="
Replace-Text "87×80=" "91×96="
Replace-Text "25×13=" "55×24="
Replace-Text "76×15=" "49×60="
Replace-Text "27×23=" "88×73="
